$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6
$ws.Range("M6").Value = 1.02
$ws.Range("N6").Value = 19

# Row 7
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 8
$ws.Range("Q7").Value = 2.3
$ws.Range("R7").Value = 1.6
$ws.Range("S7").Value = 1.5
$ws.Range("T7").Value = 2.5
$ws.Range("AF7").Value = 51
$ws.Range("AS7").Value = 201
$ws.Range("AT7").Value = 2.5

# Row 9
$ws.Range("G9").Value = 1.98
$ws.Range("I9").Value = 3.55
$ws.Range("J9").Value = 2.55
$ws.Range("K9").Value = 2.07
$ws.Range("L9").Value = 4.05
$ws.Range("P9").Value = 2.72
$ws.Range("S9").Value = 1.4
$ws.Range("T9").Value = 2.52
$ws.Range("W9").Value = 6.6
$ws.Range("Y9").Value = 8.75
$ws.Range("Z9").Value = 17
$ws.Range("AA9").Value = 17
$ws.Range("AG9").Value = 9.25
$ws.Range("AH9").Value = 18
$ws.Range("AI9").Value = 12.5
$ws.Range("AJ9").Value = 50
$ws.Range("AK9").Value = 35
$ws.Range("AO9").Value = 10
$ws.Range("AP9").Value = 19.5
$ws.Range("AQ9").Value = 37
$ws.Range("AR9").Value = 75
$ws.Range("AT9").Value = 2.47
$ws.Range("AU9").Value = 7.4
$ws.Range("AW9").Value = 5.3
$ws.Range("AY9").Value = 28
$ws.Range("BB9").Value = 400

# Row 13
$ws.Range("G13").Value = 4.5
$ws.Range("H13").Value = 3.5
$ws.Range("I13").Value = 1.8
$ws.Range("J13").Value = 5
$ws.Range("L13").Value = 2.5
$ws.Range("U13").Value = 2
$ws.Range("V13").Value = 1.73
$ws.Range("Z13").Value = 51
$ws.Range("AK13").Value = 15
$ws.Range("AO13").Value = 26
$ws.Range("AR13").Value = 126
$ws.Range("AS13").Value = 301

# Row 14
$ws.Range("G14").Value = 1.85
$ws.Range("H14").Value = 3.8
$ws.Range("I14").Value = 3.8
$ws.Range("J14").Value = 2.4
$ws.Range("L14").Value = 4
$ws.Range("O14").Value = 1.18
$ws.Range("P14").Value = 4.5
$ws.Range("Q14").Value = 1.62
$ws.Range("R14").Value = 2.25
$ws.Range("U14").Value = 1.57
$ws.Range("V14").Value = 2.25
$ws.Range("W14").Value = 9.5
$ws.Range("X14").Value = 10
$ws.Range("Z14").Value = 17
$ws.Range("AD14").Value = 7.5
$ws.Range("AE14").Value = 12
$ws.Range("AH14").Value = 21
$ws.Range("AL14").Value = 29
$ws.Range("AM14").Value = 126
$ws.Range("AO14").Value = 9.5
$ws.Range("AQ14").Value = 29
$ws.Range("AW14").Value = 6
$ws.Range("AX14").Value = 19
$ws.Range("AY14").Value = 23
$ws.Range("AZ14").Value = 51
$ws.Range("BA14").Value = 67

# Row 15
$ws.Range("G15").Value = 3.3
$ws.Range("K15").Value = 1.83
$ws.Range("O15").Value = 1.57
$ws.Range("P15").Value = 2.25
$ws.Range("Q15").Value = 2.88
$ws.Range("R15").Value = 1.4
$ws.Range("U15").Value = 2.25
$ws.Range("V15").Value = 1.57
$ws.Range("AB15").Value = 51
$ws.Range("AC15").Value = 5.5
$ws.Range("AK15").Value = 26
$ws.Range("AY15").Value = 34
$ws.Range("BA15").Value = 101
$ws.Range("BB15").Value = 351

# Row 17
$ws.Range("M17").Value = 1.14
$ws.Range("N17").Value = 5.5
$ws.Range("Q17").Value = 3.4
$ws.Range("R17").Value = 1.33
$ws.Range("T17").Value = 2.08

# Row 19
$ws.Range("G19").Value = 3.05
$ws.Range("H19").Value = 3.2
$ws.Range("I19").Value = 2.22
$ws.Range("J19").Value = 3.65
$ws.Range("K19").Value = 2.05
$ws.Range("L19").Value = 2.8
$ws.Range("O19").Value = 1.28
$ws.Range("P19").Value = 3.05
$ws.Range("R19").Value = 1.83
$ws.Range("S19").Value = 1.39
$ws.Range("T19").Value = 2.55
$ws.Range("V19").Value = 2
$ws.Range("X19").Value = 16
$ws.Range("Y19").Value = 10.75
$ws.Range("Z19").Value = 40
$ws.Range("AA19").Value = 27
$ws.Range("AC19").Value = 9.75
$ws.Range("AD19").Value = 6.2
$ws.Range("AK19").Value = 17
$ws.Range("AL19").Value = 25
$ws.Range("AN19").Value = 5
$ws.Range("AO19").Value = 17
$ws.Range("AP19").Value = 24
$ws.Range("AQ19").Value = 80
$ws.Range("AR19").Value = 120
$ws.Range("AT19").Value = 2.52
$ws.Range("AW19").Value = 4.15
$ws.Range("AX19").Value = 11.5
$ws.Range("AY19").Value = 18.5
$ws.Range("BA19").Value = 70
$ws.Range("BB19").Value = 200
